# Sleep Diary workbook update
# Fills in the "Day 4" column (column F) of the last week's diary block
# (rows 119-135, week starting 2025-1-03) with the recorded sleep data,
# mirroring the values already captured for the other six days.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 您今天早上几点醒来? (what time did you wake up this morning?) -> 07:47
$ws.Range("F122").Value = 0.32430555555555557
# 您今天几点起床? (what time did you get out of bed?) -> 08:12
$ws.Range("F123").Value = 0.34166666666666667
# 您昨晚几点上床? (what time did you go to bed last night?) -> 23:30
$ws.Range("F124").Value = 0.97916666666666663
# 您昨晚几点熄灯? (what time did you turn off the lights last night?) -> 23:30
$ws.Range("F125").Value = 0.97916666666666663
# 您昨晚熄灯后花了多长时间入睡(分钟)? (minutes to fall asleep)
$ws.Range("F126").Value = 5
# 您整晚醒来几次? (number of times woken up)
$ws.Range("F127").Value = 1
# 您整晚总共醒了多长时间(分钟)? (total minutes awake)
$ws.Range("F128").Value = 6
# 您整晚总共睡了多长时间(分钟)? (total minutes slept)
$ws.Range("F129").Value = 480
# 您昨晚睡前是否使用了影响睡眠的物质...? (substances before sleep)
$ws.Range("F130").Value = "无"
# 您昨晚睡前是否使用了电子产品...? (electronics before sleep)
$ws.Range("F131").Value = "无"
# 您昨晚睡前的身体紧张程度如何? (physical tension level)
$ws.Range("F132").Value = 3
# 您昨晚睡前的精神紧张程度如何? (mental tension level)
$ws.Range("F133").Value = 4
# 您整晚的睡眠质量如何? (sleep quality)
$ws.Range("F134").Value = 3
# 您昨天白天是否小睡? (daytime nap)
$ws.Range("F135").Value = "无"
